$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the date/time number format from the cell above (A11) so the new
# A12 cell reuses the existing style instead of creating a new one.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row of data (row 12)
$ws.Range("A12").Value = 42619.894409722219
$ws.Range("B12").Value = -8
$ws.Range("C12").Value = 54
$ws.Range("D12").Value = 44
$ws.Range("E12").Value = 54
$ws.Range("F12").Value = 76
$ws.Range("G12").Value = 7595
$ws.Range("H12").Value = 13691
$ws.Range("I12").Value = 1370
$ws.Range("J12").Value = 260
$ws.Range("K12").Value = 212
$ws.Range("L12").Value = 5
$ws.Range("M12").Value = 16
$ws.Range("N12").Value = "Bag"
